$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 7 and 8 (footprint/value text changed) ---
$ws.Range("A7").Value = "100nF Cap"
$ws.Range("A8").Value = "560 Ohm Resistor"

# --- New "NEW" section header at row 10 (bold, like the row-1 header) ---
$ws.Range("A10").Value = "NEW"
$ws.Range("A10").Font.Bold = $true

# --- Row 11: MCU (unchanged part, repeated in the new section) ---
$ws.Range("A11").Value = "MCU"
$url11 = "https://www.lcsc.com/product-detail/Microcontroller-Units-MCUs-MPUs-SOCs_Microchip-Tech-ATTINY414-SSNR_C189365.html"
$ws.Range("B11").Value = $url11
$ws.Hyperlinks.Add($ws.Range("B11"), $url11)
$ws.Range("B11").Style = "Link"

# --- Row 12: CR2032 Battery holder (unchanged part, repeated) ---
$ws.Range("A12").Value = "CR2032 Battery holder"
$url12 = "https://www.lcsc.com/product-detail/span-style-background-color-ff0-Battery-span-Connectors_Q-J-C70373_C70373.html"
$ws.Range("B12").Value = $url12
$ws.Hyperlinks.Add($ws.Range("B12"), $url12)
$ws.Range("B12").Style = "Link"

# --- Row 13: IR LED (footprint replaced -> new LCSC link) ---
$ws.Range("A13").Value = "IR LED"
$url13 = "https://www.lcsc.com/product-detail/Infrared-IR-LEDs_XSSY-XS-IR04A05-802_C5205266.html"
$ws.Range("B13").Value = $url13
$ws.Hyperlinks.Add($ws.Range("B13"), $url13)
$ws.Range("B13").Style = "Link"

# --- Row 14: LED Rot (unchanged part, repeated) ---
$ws.Range("A14").Value = "LED Rot"
$url14 = "https://www.lcsc.com/product-detail/Light-Emitting-Diodes-span-style-background-color-ff0-LED-span_XINGLIGHT-XL-1608SURC-06_C965799.html"
$ws.Range("B14").Value = $url14
$ws.Hyperlinks.Add($ws.Range("B14"), $url14)
$ws.Range("B14").Style = "Link"

# --- Row 15: Switch (unchanged part, repeated) ---
$ws.Range("A15").Value = "Switch"
$url15 = "https://www.lcsc.com/product-detail/span-style-background-color-ff0-Slide-span-Switches_SHOU-HAN-MSK12CO2-SZ_C2681568.html"
$ws.Range("B15").Value = $url15
$ws.Hyperlinks.Add($ws.Range("B15"), $url15)
$ws.Range("B15").Style = "Link"

# --- Row 16: 100nF Cap (unchanged part, repeated) ---
$ws.Range("A16").Value = "100nF Cap"
$url16 = "https://www.lcsc.com/product-detail/Multilayer-Ceramic-Capacitors-MLCC-SMD-SMT_CCTC-TCC0603X7R104K500CT_C282519.html"
$ws.Range("B16").Value = $url16
$ws.Hyperlinks.Add($ws.Range("B16"), $url16)
$ws.Range("B16").Style = "Link"

# --- Row 17: 560 Ohm Resistor (footprint replaced -> new LCSC link) ---
$ws.Range("A17").Value = "560 Ohm Resistor"
$url17 = "https://www.lcsc.com/product-detail/Chip-Resistor-Surface-Mount_UNI-ROYAL-Uniroyal-Elec-0603WAJ0561T5E_C25247.html"
$ws.Range("B17").Value = $url17
$ws.Hyperlinks.Add($ws.Range("B17"), $url17)
$ws.Range("B17").Style = "Link"

# --- Row 18: N-MOSFET (new part, plain text link - no hyperlink object) ---
$ws.Range("A18").Value = "N-MOSFET"
$ws.Range("B18").Value = "https://www.lcsc.com/product-detail/MOSFETs_YONGYUTAI-SI2302_C2891732.html"

# --- Row 19: 1M Ohm Resistor (new part, plain text link - no hyperlink object) ---
$ws.Range("A19").Value = "1M Ohm Resistor"
$ws.Range("B19").Value = "https://www.lcsc.com/product-detail/Chip-Resistor-Surface-Mount_Sunway-SC0603F1004F2BNRH_C3152128.html"

# --- Final selection matches the saved view in the target file ---
$null = $ws.Range("B17").Select()
